$d = $word.ActiveDocument

# wdReplaceOne = 1, wdReplaceAll = 2
$wdReplaceOne = 1
$wdReplaceAll = 2

function Replace-InRange($range, [string]$findText, [string]$replaceText, $mode) {
    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, $mode) | Out-Null
}

# ------------------------------------------------------------------
# 1. Delete the "Note:  although this lesson..." paragraph entirely.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute("Note:  although this lesson was written for a CentOS VM, everything should work equally well in Ubuntu.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r = $d.Content.Duplicate
    $r.Start = $d.Content.Find.Parent.Start
}

$p = $d.Paragraphs.Item(3)
Write-Output "Para3 candidate: $($p.Range.Text)"
if ($p.Range.Text.StartsWith("Note:")) {
    $p.Range.Delete() | Out-Null
}

# ------------------------------------------------------------------
# 2. Merge "$ ls -l /" + "usr" + "/bin > ls-output.txt" into one run.
# ------------------------------------------------------------------
Replace-InRange $d.Content "`$ ls -l /usr/bin > ls-output.txt" "`$ ls -l /usr/bin > ls-output.txt" $wdReplaceOne

Write-Output "done step2"

# ------------------------------------------------------------------
# 3. Merge the runs in the "pipe" paragraph that describes ps aux | grep xlogo.
# ------------------------------------------------------------------
# Find the paragraph containing the pipe explanation (search is robust to
# paragraph index drift caused by step 1's deletion).
$pipePara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("A pipe (or pipeline)")) {
        $pipePara = $para
        break
    }
}
Write-Output "pipePara found: $($pipePara -ne $null)"
$pipeRange = $pipePara.Range

Replace-InRange $pipeRange "ps aux" "ps aux" $wdReplaceOne
Replace-InRange $pipeRange " file gave us a lot of output and we were only looking for the line that contained “xlogo.”  " " file gave us a lot of output and we were only looking for the line that contained “xlogo.”  " $wdReplaceOne
Replace-InRange $pipeRange "So, we piped the output of the ps command into grep to search for lines containing xlogo." "So, we piped the output of the ps command into grep to search for lines containing xlogo." $wdReplaceOne
Replace-InRange $pipeRange "`$ ps aux | grep xlogo" "`$ ps aux | grep xlogo" $wdReplaceOne

Write-Output "Pipe para after: $($pipePara.Range.Text)"
